$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 191; this shifts the existing
# rows 191:274 down to 192:275 (row 274's old data ends up as new row 275).
$ws.Rows("191").Insert()

# Populate the newly inserted row 191 with the new data point.
$ws.Range("A191").Value = 9
$ws.Range("B191").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C191").Value = "Metropolitana"
$ws.Range("D191").NumberFormat = $ws.Range("D192").NumberFormat
$ws.Range("D191").Value = 44825
$ws.Range("E191").Value = 13
$ws.Range("F191").Value = 100112026
$ws.Range("G191").Value = "Haba"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 50
$ws.Range("K191").Value = 10000
$ws.Range("L191").Value = 10000
$ws.Range("M191").Value = 10000
$ws.Range("N191").Value = "`$/saco 25 kilos"
$ws.Range("O191").Value = "Provincia de Limarí"
$ws.Range("P191").Value = 400
$ws.Range("Q191").Value = 25
$ws.Range("R191").Value = "Hortaliza"
